$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Visual"
$ws.Range("B2").Value = "Code"
$ws.Range("C2").Value = "Git"

$ws.Range("C2").Select()
